$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update bug status (C10) from OPEN to FIXED
$ws.Range("C10").Value = "FIXED"

# Update last modified date (G10) from 40262 to 40266 (2010-03-29)
$ws.Range("G10").Value = 40266

# Adjust the window view / selection to match new state
$ws.Range("C8").Select()
